# "initial version of stimuli order" -- repopulate the per-trial sequence
# (image path / German verb / category) for rows 2-33 on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: image (col B), word (col C), category (col D)
$stimuli = @(
    @('car/car095.png', 'liefern', 'car'),
    @('dog/dog091.png', 'opfern', 'dog'),
    @('car/car064.png', 'rücken', 'car'),
    @('car/car103.png', 'kaufen', 'car'),
    @('car/car104.png', 'wenden', 'car'),
    @('car/car111.png', 'rasen', 'car'),
    @('car/car065.png', 'fliehen', 'car'),
    @('dog/dog072.png', 'biegen', 'dog'),
    @('car/car085.png', 'sondern', 'car'),
    @('car/car093.png', 'drohen', 'car'),
    @('car/car088.png', 'haken', 'car'),
    @('car/car072.png', 'enden', 'car'),
    @('dog/dog084.png', 'segeln', 'dog'),
    @('dog/dog092.png', 'klappen', 'dog'),
    @('dog/dog089.png', 'schicken', 'dog'),
    @('car/car110.png', 'füllen', 'car'),
    @('dog/dog069.png', 'stärken', 'dog'),
    @('car/car092.png', 'töten', 'car'),
    @('dog/dog110.png', 'jubeln', 'dog'),
    @('dog/dog103.png', 'fliegen', 'dog'),
    @('dog/dog088.png', 'starten', 'dog'),
    @('car/car081.png', 'regnen', 'car'),
    @('car/car066.png', 'langen', 'car'),
    @('dog/dog118.png', 'laufen', 'dog'),
    @('dog/dog085.png', 'schätzen', 'dog'),
    @('car/car073.png', 'dauern', 'car'),
    @('dog/dog101.png', 'krachen', 'dog'),
    @('car/car100.png', 'füttern', 'car'),
    @('dog/dog068.png', 'backen', 'dog'),
    @('dog/dog075.png', 'posten', 'dog'),
    @('dog/dog066.png', 'antun', 'dog'),
    @('dog/dog071.png', 'tagen', 'dog')
)

$startRow = 2
for ($i = 0; $i -lt $stimuli.Length; $i++) {
    $row = $startRow + $i
    $entry = $stimuli[$i]
    $ws.Range("B" + $row).Value = $entry[0]
    $ws.Range("C" + $row).Value = $entry[1]
    $ws.Range("D" + $row).Value = $entry[2]
}
